# Generate Report for Handback
# Refreshes the handoff/handback timestamp strings written to the report
# on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column (G2)
$overview.Range("G2").Value = "2016-08-28 15:04:29"
$overview.Range("G2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# zh-cn sheet: "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2)
$zhcn.Range("H2").Value = "2016-08-28 15:04:24"
$zhcn.Range("H2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zhcn.Range("K2").Value = "2016-08-28 15:04:50"
$zhcn.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# de-de sheet: "Correspond Handback DateTime" (K2)
$dede.Range("K2").Value = "2016-08-28 15:04:56"
$dede.Range("K2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
